$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.188.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.346.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.182'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.97'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.51%  '
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '687.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.890.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.201.86'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.357.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.07%  '
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '563.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.91%  '
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("E33").Value = '  +1.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.72%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.714.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.137'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("E40").Value = '  +1.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.39%  '
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("E45").Value = '  +1.35%  '
$ws.Range("E46").Value = '  +2.12%  '
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("E51").Value = '  -0.80%  '
